$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# Updated capital structure database
# Columns D, G-L, U-AC, AD, AF-AK, AL-AM, AN-AQ for rows 2 and 3

foreach ($r in 2,3) {
    $ws.Range("D$r").Value = -0.208

    $ws.Range("G$r").Value = -0.2736111111111111
    $ws.Range("H$r").Value = -0.2736111111111111
    $ws.Range("I$r").Value = -0.45625
    $ws.Range("J$r").Value = -0.45625
    $ws.Range("K$r").Value = -0.89
    $ws.Range("L$r").Value = -0.6180555555555556

    $ws.Range("U$r").Value = 0.292
    $ws.Range("V$r").Value = 0.2994871794871795
    $ws.Range("W$r").Value = 1.145431145431145
    $ws.Range("X$r").Value = 0.4087167192407813
    $ws.Range("Y$r").Value = 0.7367144261903642
    $ws.Range("Z$r").Value = 0.2912032355915066
    $ws.Range("AA$r").Value = -0.1328614762386249
    $ws.Range("AB$r").Value = 0.09087079008916338
    $ws.Range("AC$r").Value = -0.2237322663277883

    $ws.Range("AD$r").Value = 6.28
    $ws.Range("AF$r").Value = 6.28
    $ws.Range("AG$r").Value = 5.988
    $ws.Range("AH$r").Value = 0.8656099241902137
    $ws.Range("AI$r").Value = 1.37417943107221
    $ws.Range("AJ$r").Value = 0.8599741490736752
    $ws.Range("AK$r").Value = 1.399719495091164

    $ws.Range("AL$r").Value = 0.107
    $ws.Range("AM$r").Value = 0.107

    $ws.Range("AN$r").Value = -12.41106719367589
    $ws.Range("AO$r").Value = -6.140186915887851
    $ws.Range("AP$r").Value = -11.83399209486166
    $ws.Range("AQ$r").Value = -6.140186915887851
}
